# Estado de Cuenta - actualizacion de base de datos (NIT-8000551164)
# - Elimina EC anteriores y se agregan nuevos (nuevas filas de deudores/periodos)
# - Se modifica base de datos (totales, cantidad de trabajadores/periodos)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the existing worker name text verbatim (keeps the mojibake byte-for-byte
# identical instead of retyping it) before we start shuffling rows around.
$sandraName = $ws.Range("D16").Value()

# ------------------------------------------------------------------
# 1) Make room: insert 5 new data rows right above the last (bottom
#    bordered) row of the table, then copy the formatting of the row
#    directly above down into the freshly inserted rows.
# ------------------------------------------------------------------
$ws.Rows("25:29").Insert()

$ws.Range("B24:J24").Copy()
$ws.Range("B25:J29").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# ------------------------------------------------------------------
# 2) Header / summary cells
# ------------------------------------------------------------------
$ws.Range("E11").Value = 2508940      # VALOR MORA
$ws.Range("C13").Value = 3            # Cant. Trabajadores
$ws.Range("F13").Value = 11           # Cant. Periodos

# ------------------------------------------------------------------
# 3) Data table (rows 16-30) - new debtor / period detail
# ------------------------------------------------------------------
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "45706252"
$ws.Range("D16").Value = "YANERIS MARGARITA MENDOZA GUERRA"
$ws.Range("E16").Value = "2505"
$ws.Range("F16").Value = 84000
$ws.Range("G16").Value = 908526

$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "45706252"
$ws.Range("D17").Value = "YANERIS MARGARITA MENDOZA GUERRA"
$ws.Range("E17").Value = "2504"
$ws.Range("F17").Value = 84000
$ws.Range("G17").Value = 908526

$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "45706252"
$ws.Range("D18").Value = "YANERIS MARGARITA MENDOZA GUERRA"
$ws.Range("E18").Value = "2503"
$ws.Range("F18").Value = 84000
$ws.Range("G18").Value = 908526

$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "1044918760"
$ws.Range("D19").Value = "JAIRO JOSE HERNANDEZ BELLO"
$ws.Range("E19").Value = "2507"
$ws.Range("F19").Value = 56940
$ws.Range("G19").Value = 1423500

$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "1098623598"
$ws.Range("D20").Value = $sandraName
$ws.Range("E20").Value = "2507"
$ws.Range("F20").Value = 200000
$ws.Range("G20").Value = 5000000

$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "1098623598"
$ws.Range("D21").Value = $sandraName
$ws.Range("E21").Value = "2506"
$ws.Range("F21").Value = 200000
$ws.Range("G21").Value = 5000000

$ws.Range("B22").Value = "CC"
$ws.Range("C22").Value = "1098623598"
$ws.Range("D22").Value = $sandraName
$ws.Range("E22").Value = "2505"
$ws.Range("F22").Value = 200000
$ws.Range("G22").Value = 5000000

$ws.Range("B23").Value = "CC"
$ws.Range("C23").Value = "1098623598"
$ws.Range("D23").Value = $sandraName
$ws.Range("E23").Value = "2504"
$ws.Range("F23").Value = 200000
$ws.Range("G23").Value = 5000000

$ws.Range("B24").Value = "CC"
$ws.Range("C24").Value = "1098623598"
$ws.Range("D24").Value = $sandraName
$ws.Range("E24").Value = "2503"
$ws.Range("F24").Value = 200000
$ws.Range("G24").Value = 5000000

$ws.Range("B25").Value = "CC"
$ws.Range("C25").Value = "1098623598"
$ws.Range("D25").Value = $sandraName
$ws.Range("E25").Value = "2502"
$ws.Range("F25").Value = 200000
$ws.Range("G25").Value = 5000000

$ws.Range("B26").Value = "CC"
$ws.Range("C26").Value = "1098623598"
$ws.Range("D26").Value = $sandraName
$ws.Range("E26").Value = "2501"
$ws.Range("F26").Value = 200000
$ws.Range("G26").Value = 5000000

$ws.Range("B27").Value = "CC"
$ws.Range("C27").Value = "1098623598"
$ws.Range("D27").Value = $sandraName
$ws.Range("E27").Value = "2412"
$ws.Range("F27").Value = 200000
$ws.Range("G27").Value = 5000000

$ws.Range("B28").Value = "CC"
$ws.Range("C28").Value = "1098623598"
$ws.Range("D28").Value = $sandraName
$ws.Range("E28").Value = "2411"
$ws.Range("F28").Value = 200000
$ws.Range("G28").Value = 5000000

$ws.Range("B29").Value = "CC"
$ws.Range("C29").Value = "1098623598"
$ws.Range("D29").Value = $sandraName
$ws.Range("E29").Value = "2410"
$ws.Range("F29").Value = 200000
$ws.Range("G29").Value = 5000000

$ws.Range("B30").Value = "CC"
$ws.Range("C30").Value = "1098623598"
$ws.Range("D30").Value = $sandraName
$ws.Range("E30").Value = "2409"
$ws.Range("F30").Value = 200000
$ws.Range("G30").Value = 5000000

# ------------------------------------------------------------------
# 4) Column D autofit (the new, longer worker name widens the column)
# ------------------------------------------------------------------
$ws.Columns("D").AutoFit()

Write-Output ("Final used range: " + $ws.UsedRange.Address())
